$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "<rpc-reply message-id=""urn:uuid:b553d70a-307c-44a1-859a-5c44db234b87"">`n  <data>`n    <network-instances>`n      <network-instance>`n        <name>Prueba_LxVPN</name>`n        <config>`n          <name>Prueba_LxVPN</name>`n          <type>oc-ni-types:L3VRF</type>`n        </config>`n        <interfaces>`n          <interface>`n            <id>GigabitEthernet0/3/2</id>`n            <config>`n              <id>GigabitEthernet0/3/2</id>`n              <interface>GigabitEthernet0/3/2</interface>`n              <subinterface>0</subinterface>`n            </config>`n          </interface>`n        </interfaces>`n        <protocols>`n          <protocol>`n            <identifier>oc-pol-types:OSPF</identifier>`n            <name>22</name>`n            <config>`n              <identifier>oc-pol-types:OSPF</identifier>`n              <name>22</name>`n            </config>`n            <ospfv2>`n              <global>`n                <config>`n                  <router-id>172.16.1.3</router-id>`n                </config>`n              </global>`n            </ospfv2>`n          </protocol>`n          <protocol>`n            <identifier>oc-pol-types:STATIC</identifier>`n            <name>default</name>`n            <config>`n              <identifier>oc-pol-types:STATIC</identifier>`n              <name>default</name>`n            </config>`n          </protocol>`n          <protocol>`n            <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>`n            <name>default</name>`n            <config>`n              <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>`n              <name>default</name>`n            </config>`n          </protocol>`n        </protocols>`n      </network-instance>`n    </network-instances>`n  </data>`n</rpc-reply>`n"

$ws.Range("F2").Value = $newValue
